# "add thank you slide"
#
# Appends a new slide (the 3rd slide) using the "Title Only" layout and
# sets its title placeholder to "Thank You", positioned/sized to match
# the target deck. PowerPoint stores Shape.Left/Top/Width/Height in
# points (1 pt = 12700 EMU); the literals below are the point values
# whose single-precision round trip lands exactly on the target EMU
# offsets (4628965, 2451377) and extents (3645023, 1325563).

$p = $ppt.ActivePresentation

$slide = $p.Slides.Add($p.Slides.Count + 1, 11)  # 11 = ppLayoutTitleOnly

$title = $slide.Shapes.Item(1)
$title.Name = "제목 1"
$title.Left = 364.4854431152344
$title.Top = 193.02182006835938
$title.Width = 287.00970458984375
$title.Height = 104.37504577636719
$title.TextFrame.TextRange.Text = "Thank You"
